$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "2025-08-11 09:22:29"
$ws.Range("B4").Value = "hsenbyomi@gmail.com"
$ws.Range("C4").Value = "password123"
$ws.Range("D4").Value = "http://localhost:3000/uploads/1754904149624-undefined"
